# Updates cryptos list values (price & 1h volume change) to match the
# latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.271.08"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "1.775.52"
$ws.Range("E3").Value = "  +3.60%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'313.31"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "'0.5207"
$ws.Range("E7").Value = "  +8.76%  "
$ws.Range("D8").Value = "'0.3696"
$ws.Range("E8").Value = "  +7.14%  "
$ws.Range("D9").Value = "'42.79"
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("D10").Value = "'0.07394"
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("E11").Value = "  +4.56%  "
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Value = "'20.50"
$ws.Range("E13").Value = "  +3.32%  "
$ws.Range("E14").Value = "  +3.71%  "
$ws.Range("D15").Value = "1.766.84"
$ws.Range("E15").Value = "  +3.12%  "
$ws.Range("D16").Value = "'6.965"
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("D17").Value = "'89.02"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "'0.00001049"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").Value = "'0.06444"
$ws.Range("E19").Value = "  +1.27%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("E21").Value = "  +1.61%  "
$ws.Range("D22").Value = "'5.819"
$ws.Range("E22").Value = "  +3.68%  "
$ws.Range("D23").Value = "27.300.64"
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("E24").Value = "  +4.02%  "
$ws.Range("D25").Value = "'2.122"
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("D26").Value = "'155.10"
$ws.Range("E26").Value = "  +2.08%  "
$ws.Range("D27").Value = "'20.20"
$ws.Range("E27").Value = "  +2.81%  "
$ws.Range("D28").Value = "1.976.60"
$ws.Range("E28").Value = "  +3.53%  "
$ws.Range("D29").Value = "'2.325"
$ws.Range("D30").Value = "'121.28"
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("E31").Value = "  +4.52%  "
$ws.Range("D32").Value = "'0.09795"
$ws.Range("E32").Value = "  +5.56%  "
$ws.Range("D33").Value = "'5.572"
$ws.Range("E33").Value = "  +4.99%  "
$ws.Range("D34").Value = "'3.621"
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("D35").Value = "'0.02242"
$ws.Range("E35").Value = "  +1.90%  "
$ws.Range("D36").Value = "'0.05977"
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").Value = "'11.25"
$ws.Range("E37").Value = "  +1.80%  "
$ws.Range("D38").Value = "'4.840"
$ws.Range("E38").Value = "  +2.24%  "
$ws.Range("D39").Value = "'0.6139"
$ws.Range("E39").Value = "  +3.62%  "
$ws.Range("D40").Value = "'0.2019"
$ws.Range("E40").Value = "  +0.83%  "
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("D42").Value = "'8.092"
$ws.Range("E42").Value = "  +8.36%  "
$ws.Range("D43").Value = "'1.141"
$ws.Range("E43").Value = "  +2.62%  "
$ws.Range("E44").Value = "  +2.79%  "
$ws.Range("E45").Value = "  +2.82%  "
$ws.Range("D46").Value = "'3.628"
$ws.Range("E46").Value = "  +1.61%  "
$ws.Range("D47").Value = "'121.26"
$ws.Range("E47").Value = "  +2.16%  "
$ws.Range("D48").Value = "'1.886"
$ws.Range("E48").Value = "  +2.72%  "
$ws.Range("D49").Value = "'1.116"
$ws.Range("E49").Value = "  +3.01%  "
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("D51").Value = "'70.56"
$ws.Range("E51").Value = "  +1.34%  "
